$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 2559.8  # H51: was 3097.7778
$ws.Cells.Item(51, 9).Value = 2050  # I51: was 3945
$ws.Cells.Item(51, 10).Value = 3324.5  # J51: was 2420
$ws.Cells.Item(51, 11).Value = 2050  # K51: was 3945
$ws.Cells.Item(51, 12).Value = 3324.5  # L51: was 2420
$ws.Cells.Item(51, 13).Value = -1566  # M51: was -3461
$ws.Cells.Item(51, 14).Value = -4292.5  # N51: was -3388

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 51666.668  # H134: was 55000
$ws.Cells.Item(134, 10).Value = 51666.668  # J134: was 55000
$ws.Cells.Item(134, 12).Value = 51666.668  # L134: was 55000
$ws.Cells.Item(134, 14).Value = -61806.668  # N134: was -65140

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 55557460  # H137: was 62501836
$ws.Cells.Item(137, 9).Value = 62501204  # I137: was 66667956
$ws.Cells.Item(137, 10).Value = 7501.5  # J137: was 10003
$ws.Cells.Item(137, 11).Value = 187503612  # K137: was 200003868
$ws.Cells.Item(137, 12).Value = 22504.5  # L137: was 30009
$ws.Cells.Item(137, 13).Value = -187501062  # M137: was -200001318
$ws.Cells.Item(137, 14).Value = -27604.5  # N137: was -35109

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13656.036  # H32: was 15210.284
$ws.Cells.Item(32, 9).Value = 1564.4684  # I32: was 1478.942
$ws.Cells.Item(32, 11).Value = 1564.4684  # K32: was 1478.942
$ws.Cells.Item(32, 13).Value = -1277.4684  # M32: was -1191.942

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 1037  # H35: was 0
$ws.Cells.Item(35, 9).Value = 1037  # I35: was 0
$ws.Cells.Item(35, 11).Value = 1037  # K35: was 0
$ws.Cells.Item(35, 13).Value = -631  # M35: was None

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 6914.923  # H74: was 8562.936
$ws.Cells.Item(74, 9).Value = 1186.5  # I74: was 1451.9524
$ws.Cells.Item(74, 10).Value = 21496.363  # J74: was 23496
$ws.Cells.Item(74, 11).Value = 1186.5  # K74: was 1451.9524
$ws.Cells.Item(74, 12).Value = 21496.363  # L74: was 23496
$ws.Cells.Item(74, 13).Value = -312.5  # M74: was -577.9523999999999
$ws.Cells.Item(74, 14).Value = -23244.363  # N74: was -25244

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 6914.923  # H77: was 8562.936
$ws.Cells.Item(77, 9).Value = 1186.5  # I77: was 1451.9524
$ws.Cells.Item(77, 10).Value = 21496.363  # J77: was 23496
$ws.Cells.Item(77, 11).Value = 5932.5  # K77: was 7259.762
$ws.Cells.Item(77, 12).Value = 107481.815  # L77: was 117480
$ws.Cells.Item(77, 13).Value = -1564.5  # M77: was -2891.762
$ws.Cells.Item(77, 14).Value = -116217.815  # N77: was -126216

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 5266.5  # H88: was 10950
$ws.Cells.Item(88, 9).Value = 1899.75  # I88: was 1900
$ws.Cells.Item(88, 10).Value = 12000  # J88: was 20000
$ws.Cells.Item(88, 11).Value = 1899.75  # K88: was 1900
$ws.Cells.Item(88, 12).Value = 12000  # L88: was 20000
$ws.Cells.Item(88, 13).Value = -1493.75  # M88: was -1494
$ws.Cells.Item(88, 14).Value = -12812  # N88: was -20812

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 5266.5  # H91: was 10950
$ws.Cells.Item(91, 9).Value = 1899.75  # I91: was 1900
$ws.Cells.Item(91, 10).Value = 12000  # J91: was 20000
$ws.Cells.Item(91, 11).Value = 1899.75  # K91: was 1900
$ws.Cells.Item(91, 12).Value = 12000  # L91: was 20000
$ws.Cells.Item(91, 13).Value = -495.75  # M91: was -496
$ws.Cells.Item(91, 14).Value = -14808  # N91: was -22808

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1642  # H102: was 1635
$ws.Cells.Item(102, 9).Value = 1642  # I102: was 1635
$ws.Cells.Item(102, 11).Value = 1642  # K102: was 1635
$ws.Cells.Item(102, 13).Value = -20  # M102: was -13

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 453.2143  # H110: was 469.5926
$ws.Cells.Item(110, 9).Value = 428.75  # I110: was 446.91306
$ws.Cells.Item(110, 11).Value = 428.75  # K110: was 446.91306
$ws.Cells.Item(110, 13).Value = 1616.25  # M110: was 1598.08694

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1119.5483  # H20: was 1463.6
$ws.Cells.Item(20, 9).Value = 1124  # I20: was 1514.3334
$ws.Cells.Item(20, 10).Value = 1112.5  # J20: was 1387.5
$ws.Cells.Item(20, 11).Value = 1124  # K20: was 1514.3334
$ws.Cells.Item(20, 12).Value = 1112.5  # L20: was 1387.5
$ws.Cells.Item(20, 13).Value = -877  # M20: was -1267.3334
$ws.Cells.Item(20, 14).Value = -1606.5  # N20: was -1881.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(46, 8).Value = 5300  # H46: was 0
$ws.Cells.Item(46, 10).Value = 5300  # J46: was 0
$ws.Cells.Item(46, 12).Value = 5300  # L46: was 0
$ws.Cells.Item(46, 14).Value = -5896  # N46: was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6653.8  # H86: was 5898.478
$ws.Cells.Item(86, 9).Value = 1858.6666  # I86: was 1856
$ws.Cells.Item(86, 10).Value = 13846.5  # J86: was 9604.083000000001
$ws.Cells.Item(86, 11).Value = 1858.6666  # K86: was 1856
$ws.Cells.Item(86, 12).Value = 13846.5  # L86: was 9604.083000000001
$ws.Cells.Item(86, 13).Value = -735.6666  # M86: was -733
$ws.Cells.Item(86, 14).Value = -16092.5  # N86: was -11850.083

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 6653.8  # H89: was 5898.478
$ws.Cells.Item(89, 9).Value = 1858.6666  # I89: was 1856
$ws.Cells.Item(89, 10).Value = 13846.5  # J89: was 9604.083000000001
$ws.Cells.Item(89, 11).Value = 9293.333000000001  # K89: was 9280
$ws.Cells.Item(89, 12).Value = 69232.5  # L89: was 48020.415
$ws.Cells.Item(89, 13).Value = -3677.333000000001  # M89: was -3664
$ws.Cells.Item(89, 14).Value = -80464.5  # N89: was -59252.415

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1103  # H94: was 1109
$ws.Cells.Item(94, 9).Value = 1259.1818  # I94: was 1220.9166
$ws.Cells.Item(94, 10).Value = 759.4  # J94: was 773.25
$ws.Cells.Item(94, 11).Value = 1259.1818  # K94: was 1220.9166
$ws.Cells.Item(94, 12).Value = 759.4  # L94: was 773.25
$ws.Cells.Item(94, 13).Value = -808.1818000000001  # M94: was -769.9166
$ws.Cells.Item(94, 14).Value = -1661.4  # N94: was -1675.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2982.7354  # H105: was 3085.0938
$ws.Cells.Item(105, 9).Value = 2719  # I105: was 2884.2917
$ws.Cells.Item(105, 10).Value = 4000  # J105: was 3687.5
$ws.Cells.Item(105, 11).Value = 2719  # K105: was 2884.2917
$ws.Cells.Item(105, 12).Value = 4000  # L105: was 3687.5
$ws.Cells.Item(105, 13).Value = -972  # M105: was -1137.2917
$ws.Cells.Item(105, 14).Value = -7494  # N105: was -7181.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1039.1538  # H107: was 999.4375
$ws.Cells.Item(107, 9).Value = 685.8570999999999  # I107: was 643.3333
$ws.Cells.Item(107, 10).Value = 1451.3334  # J107: was 1457.2858
$ws.Cells.Item(107, 11).Value = 685.8570999999999  # K107: was 643.3333
$ws.Cells.Item(107, 12).Value = 1451.3334  # L107: was 1457.2858
$ws.Cells.Item(107, 13).Value = 1234.1429  # M107: was 1276.6667
$ws.Cells.Item(107, 14).Value = -5291.3334  # N107: was -5297.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2673.1924  # H134: was 2521.5356
$ws.Cells.Item(134, 9).Value = 2248  # I134: was 2146.2693
$ws.Cells.Item(134, 10).Value = 5933  # J134: was 7400
$ws.Cells.Item(134, 11).Value = 6744  # K134: was 6438.8079
$ws.Cells.Item(134, 12).Value = 17799  # L134: was 22200
$ws.Cells.Item(134, 13).Value = -4209  # M134: was -3903.8079
$ws.Cells.Item(134, 14).Value = -22869  # N134: was -27270

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 83333.336  # H23: was 26995
$ws.Cells.Item(23, 9).Value = 50000  # I23: was 26995
$ws.Cells.Item(23, 10).Value = 100000  # J23: was 0
$ws.Cells.Item(23, 11).Value = 50000  # K23: was 26995
$ws.Cells.Item(23, 12).Value = 100000  # L23: was 0
$ws.Cells.Item(23, 13).Value = -49760  # M23: was -26755
$ws.Cells.Item(23, 14).Value = -100480  # N23: was None

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(27, 8).Value = 83333.336  # H27: was 26995
$ws.Cells.Item(27, 9).Value = 50000  # I27: was 26995
$ws.Cells.Item(27, 10).Value = 100000  # J27: was 0
$ws.Cells.Item(27, 11).Value = 50000  # K27: was 26995
$ws.Cells.Item(27, 12).Value = 100000  # L27: was 0
$ws.Cells.Item(27, 13).Value = -49808  # M27: was -26803
$ws.Cells.Item(27, 14).Value = -100384  # N27: was None

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1182.1072  # H31: was 1149.963
$ws.Cells.Item(31, 9).Value = 943.96  # I31: was 961.96
$ws.Cells.Item(31, 10).Value = 3166.6667  # J31: was 3500
$ws.Cells.Item(31, 11).Value = 943.96  # K31: was 961.96
$ws.Cells.Item(31, 12).Value = 3166.6667  # L31: was 3500
$ws.Cells.Item(31, 13).Value = -648.96  # M31: was -666.96
$ws.Cells.Item(31, 14).Value = -3756.6667  # N31: was -4090

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1182.1072  # H34: was 1149.963
$ws.Cells.Item(34, 9).Value = 943.96  # I34: was 961.96
$ws.Cells.Item(34, 10).Value = 3166.6667  # J34: was 3500
$ws.Cells.Item(34, 11).Value = 943.96  # K34: was 961.96
$ws.Cells.Item(34, 12).Value = 3166.6667  # L34: was 3500
$ws.Cells.Item(34, 13).Value = -741.96  # M34: was -759.96
$ws.Cells.Item(34, 14).Value = -3570.6667  # N34: was -3904

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 35713.43  # H62: was 35713.285
$ws.Cells.Item(62, 10).Value = 9998.5  # J62: was 9998.25
$ws.Cells.Item(62, 12).Value = 9998.5  # L62: was 9998.25
$ws.Cells.Item(62, 14).Value = -11246.5  # N62: was -11246.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 35713.43  # H65: was 35713.285
$ws.Cells.Item(65, 10).Value = 9998.5  # J65: was 9998.25
$ws.Cells.Item(65, 12).Value = 49992.5  # L65: was 49991.25
$ws.Cells.Item(65, 14).Value = -56232.5  # N65: was -56231.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1462.1428  # H94: was 1507.65
$ws.Cells.Item(94, 9).Value = 930.55554  # I94: was 977.875
$ws.Cells.Item(94, 11).Value = 930.55554  # K94: was 977.875
$ws.Cells.Item(94, 13).Value = -479.55554  # M94: was -526.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(133, 8).Value = 30745.555  # H133: was 15844.667
$ws.Cells.Item(133, 9).Value = 0  # I133: was 20296
$ws.Cells.Item(133, 10).Value = 30745.555  # J133: was 15288.25
$ws.Cells.Item(133, 11).Value = 0  # K133: was 20296
$ws.Cells.Item(133, 12).Value = 30745.555  # L133: was 15288.25
$ws.Cells.Item(133, 13).ClearContents()  # M133: was -17766
$ws.Cells.Item(133, 14).Value = -35805.555  # N133: was -20348.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3403.5833  # H134: was 3086.963
$ws.Cells.Item(134, 9).Value = 1853.5714  # I134: was 1624.2354
$ws.Cells.Item(134, 11).Value = 5560.7142  # K134: was 4872.706200000001
$ws.Cells.Item(134, 13).Value = -3025.7142  # M134: was -2337.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(135, 8).Value = 30581.5  # H135: was 42371.25
$ws.Cells.Item(135, 10).Value = 30581.5  # J135: was 42371.25
$ws.Cells.Item(135, 12).Value = 30581.5  # L135: was 42371.25
$ws.Cells.Item(135, 14).Value = -40721.5  # N135: was -52511.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(137, 8).Value = 42538  # H137: was 35837.8
$ws.Cells.Item(137, 9).Value = 10000  # I137: was 9854.5
$ws.Cells.Item(137, 10).Value = 46153.332  # J137: was 53160
$ws.Cells.Item(137, 11).Value = 10000  # K137: was 9854.5
$ws.Cells.Item(137, 12).Value = 46153.332  # L137: was 53160
$ws.Cells.Item(137, 13).Value = -4900  # M137: was -4754.5
$ws.Cells.Item(137, 14).Value = -56353.332  # N137: was -63360

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 48186.668  # H140: was 55000
$ws.Cells.Item(140, 10).Value = 48186.668  # J140: was 55000
$ws.Cells.Item(140, 12).Value = 48186.668  # L140: was 55000
$ws.Cells.Item(140, 14).Value = -58546.668  # N140: was -65360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 766.6667  # H59: was 1750
$ws.Cells.Item(59, 9).Value = 766.6667  # I59: was 2000
$ws.Cells.Item(59, 10).Value = 0  # J59: was 1500
$ws.Cells.Item(59, 11).Value = 2300.0001  # K59: was 6000
$ws.Cells.Item(59, 12).Value = 0  # L59: was 4500
$ws.Cells.Item(59, 13).Value = -1760.0001  # M59: was -5460
$ws.Cells.Item(59, 14).ClearContents()  # N59: was -5580

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1171.591  # H132: was 1170.238
$ws.Cells.Item(132, 10).Value = 1379.5  # J132: was 1393.3077
$ws.Cells.Item(132, 12).Value = 12415.5  # L132: was 12539.7693
$ws.Cells.Item(132, 14).Value = -17475.5  # N132: was -17599.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5517  # H70: was 5366.4
$ws.Cells.Item(70, 9).Value = 5440.2354  # I70: was 5458.6177
$ws.Cells.Item(70, 10).Value = 6039  # J70: was 5081.364
$ws.Cells.Item(70, 11).Value = 5440.2354  # K70: was 5458.6177
$ws.Cells.Item(70, 12).Value = 6039  # L70: was 5081.364
$ws.Cells.Item(70, 13).Value = -5170.2354  # M70: was -5188.6177
$ws.Cells.Item(70, 14).Value = -6579  # N70: was -5621.364

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5517  # H73: was 5366.4
$ws.Cells.Item(73, 9).Value = 5440.2354  # I73: was 5458.6177
$ws.Cells.Item(73, 10).Value = 6039  # J73: was 5081.364
$ws.Cells.Item(73, 11).Value = 5440.2354  # K73: was 5458.6177
$ws.Cells.Item(73, 12).Value = 6039  # L73: was 5081.364
$ws.Cells.Item(73, 13).Value = -4504.2354  # M73: was -4522.6177
$ws.Cells.Item(73, 14).Value = -7911  # N73: was -6953.364

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1681.2222  # H102: was 1687.3334
$ws.Cells.Item(102, 9).Value = 1422.625  # I102: was 1429.5
$ws.Cells.Item(102, 11).Value = 1422.625  # K102: was 1429.5
$ws.Cells.Item(102, 13).Value = 199.375  # M102: was 192.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2978.7908  # H132: was 3216.5
$ws.Cells.Item(132, 9).Value = 2698.647  # I132: was 2918.1
$ws.Cells.Item(132, 10).Value = 4037.111  # J132: was 4335.5
$ws.Cells.Item(132, 11).Value = 8095.941  # K132: was 8754.299999999999
$ws.Cells.Item(132, 12).Value = 12111.333  # L132: was 13006.5
$ws.Cells.Item(132, 13).Value = -5565.941  # M132: was -6224.299999999999
$ws.Cells.Item(132, 14).Value = -17171.333  # N132: was -18066.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 16367.5  # H133: was 18312.223
$ws.Cells.Item(133, 10).Value = 16367.5  # J133: was 18312.223
$ws.Cells.Item(133, 12).Value = 16367.5  # L133: was 18312.223
$ws.Cells.Item(133, 14).Value = -26487.5  # N133: was -28432.223

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 142891710  # H135: was 1000000000
$ws.Cells.Item(135, 10).Value = 142891710  # J135: was 1000000000
$ws.Cells.Item(135, 12).Value = 142891710  # L135: was 1000000000
$ws.Cells.Item(135, 14).Value = -142901850  # N135: was -1000010140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 61280  # H138: was 61680
$ws.Cells.Item(138, 10).Value = 61280  # J138: was 61680
$ws.Cells.Item(138, 12).Value = 61280  # L138: was 61680
$ws.Cells.Item(138, 14).Value = -71560  # N138: was -71960

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(141, 8).Value = 58485.6  # H141: was 59045.8
$ws.Cells.Item(141, 10).Value = 58485.6  # J141: was 59045.8
$ws.Cells.Item(141, 12).Value = 58485.6  # L141: was 59045.8
$ws.Cells.Item(141, 14).Value = -68845.60000000001  # N141: was -69405.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 7000  # H5: was 6250
$ws.Cells.Item(5, 9).Value = 10000  # I5: was 7500
$ws.Cells.Item(5, 10).Value = 4000  # J5: was 5000
$ws.Cells.Item(5, 11).Value = 10000  # K5: was 7500
$ws.Cells.Item(5, 12).Value = 4000  # L5: was 5000
$ws.Cells.Item(5, 13).Value = -9887  # M5: was -7387
$ws.Cells.Item(5, 14).Value = -4226  # N5: was -5226

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 3000  # H21: was 0
$ws.Cells.Item(21, 10).Value = 3000  # J21: was 0
$ws.Cells.Item(21, 12).Value = 3000  # L21: was 0
$ws.Cells.Item(21, 14).Value = -3470  # N21: was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 2999.1428  # H24: was 3000
$ws.Cells.Item(24, 10).Value = 2999.1428  # J24: was 3000
$ws.Cells.Item(24, 12).Value = 2999.1428  # L24: was 3000
$ws.Cells.Item(24, 14).Value = -3459.1428  # N24: was -3460

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 28750  # H28: was 50000
$ws.Cells.Item(28, 10).Value = 7500  # J28: was 0
$ws.Cells.Item(28, 12).Value = 7500  # L28: was 0
$ws.Cells.Item(28, 14).Value = -8196  # N28: was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 7000  # H30: was 7500
$ws.Cells.Item(30, 9).Value = 7000  # I30: was 0
$ws.Cells.Item(30, 10).Value = 0  # J30: was 7500
$ws.Cells.Item(30, 11).Value = 7000  # K30: was 0
$ws.Cells.Item(30, 12).Value = 0  # L30: was 7500
$ws.Cells.Item(30, 13).Value = -6893  # M30: was None
$ws.Cells.Item(30, 14).ClearContents()  # N30: was -7714

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(35, 8).Value = 3000  # H35: was 0
$ws.Cells.Item(35, 10).Value = 3000  # J35: was 0
$ws.Cells.Item(35, 12).Value = 3000  # L35: was 0
$ws.Cells.Item(35, 14).Value = -3580  # N35: was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 777.58826  # H113: was 800
$ws.Cells.Item(113, 9).Value = 666.75  # I113: was 683.1539
$ws.Cells.Item(113, 10).Value = 1043.6  # J113: was 1179.75
$ws.Cells.Item(113, 11).Value = 2000.25  # K113: was 2049.4617
$ws.Cells.Item(113, 12).Value = 3130.8  # L113: was 3539.25
$ws.Cells.Item(113, 13).Value = 169.75  # M113: was 120.5383000000002
$ws.Cells.Item(113, 14).Value = -7470.799999999999  # N113: was -7879.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(128, 8).Value = 49900  # H128: was 49828.57
$ws.Cells.Item(128, 10).Value = 49900  # J128: was 49828.57
$ws.Cells.Item(128, 12).Value = 49900  # L128: was 49828.57
$ws.Cells.Item(128, 14).Value = -59860  # N128: was -59788.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 13161462  # H132: was 14289533
$ws.Cells.Item(132, 9).Value = 23814092  # I132: was 27783000
$ws.Cells.Item(132, 11).Value = 71442276  # K132: was 83349000
$ws.Cells.Item(132, 13).Value = -71439746  # M132: was -83346470
